$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-06 Monday" "2025-10-07 Tuesday"

Replace-Text "99×20=" "77×14="
Replace-Text "65×28=" "57×76="
Replace-Text "52×99=" "54×50="
Replace-Text "57×47=" "30×98="
Replace-Text "12×76=" "93×73="
Replace-Text "44×49=" "14×26="
Replace-Text "93×53=" "23×75="
Replace-Text "19×65=" "87×89="
Replace-Text "11×46=" "57×15="
Replace-Text "53×14=" "59×16="
Replace-Text "80×84=" "64×44="
Replace-Text "60×53=" "27×81="
Replace-Text "44×53=" "70×32="
Replace-Text "56×24=" "93×30="
Replace-Text "52×90=" "33×31="
Replace-Text "78×23=" "85×53="
Replace-Text "52×19=" "44×40="
Replace-Text "13×16=" "35×79="
Replace-Text "26×58=" "38×74="
Replace-Text "24×32=" "41×73="
Replace-Text "50×63=" "69×28="
Replace-Text "95×81=" "92×43="
Replace-Text "49×31=" "32×83="
Replace-Text "46×66=" "53×59="
Replace-Text "49×43=" "20×46="
